$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "30.288.79"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +2.05%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.095.50"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  -0.64%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "342.47"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.22%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.51%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5292"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +2.31%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.4371"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.18%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "55.07"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +3.23%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.09334"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.35%  "
$ws.Range("E11").Value = "  +0.68%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "24.67"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.40%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "8.538"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +4.92%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.865"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "2.045.26"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.57%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "100.97"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.83%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.00001154"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.25%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "1.004"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.48%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "21.09"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.38%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.06731"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.87%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.53%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "30.274.83"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.84%  "
$ws.Range("E24").Value = "  -1.10%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.321"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.69%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "6.965"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +9.15%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "21.77"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.54%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "162.58"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.57%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.503"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.54%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "133.56"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.02%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.128"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.1052"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.663"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.00%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "6.232"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.72%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "3.916"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.92%  "
$ws.Range("E36").Value = "  -3.27%  "
$ws.Range("E37").Value = "  +1.65%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.06752"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.62%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "12.59"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.11%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.6962"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.26%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.342"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.18%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.2206"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.48%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.6780"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.08%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "14.32"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("E46").Value = "  -0.17%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.298"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +7.81%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "3.641"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.66%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.00000000349"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.87%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.212"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +5.06%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.211"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.29%  "
